# Insert a new data row at row 460 (pushing existing rows 460..571 down to 461..572)
# and populate it with a new "Ajo" (garlic) price record for Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 460, shifting everything below it down by one.
$ws.Rows.Item(460).Insert()

# Fill in the values for the newly inserted row 460.
$ws.Range("A460").Value = 3
$ws.Range("B460").Value = "Femacal de La Calera"
$ws.Range("C460").Value = "Coquimbo"
$ws.Range("D460").Value = 44855
$ws.Range("E460").Value = 5
$ws.Range("F460").Value = 100112003
$ws.Range("G460").Value = "Ajo"
$ws.Range("H460").Value = "Chino"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 105
$ws.Range("K460").Value = 14000
$ws.Range("L460").Value = 14500
$ws.Range("M460").Value = 14262
$ws.Range("N460").Value = "$/caja 10 kilos"
$ws.Range("O460").Value = "China"
$ws.Range("P460").Value = 1426
$ws.Range("Q460").Value = 10
$ws.Range("R460").Value = "Hortaliza"
